$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "214.38") are not auto-converted to numbers, matching the
# original inline-string cell type. Style is reset to Normal right
# after so no stray number-format style is left on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.717.50'
$ws.Range("E2").Value = '  -0.27%  '

$ws.Range("D3").Value = '1.629.81'
$ws.Range("E3").Value = '  -0.22%  '

$ws.Range("E4").Value = '  -0.75%  '

$ws.Range("D5").Value = '214.38'
$ws.Range("E5").Value = '  -0.53%  '

$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("E7").Value = '  -0.28%  '

$ws.Range("D8").Value = '0.255'

$ws.Range("D9").Value = '0.0632'
$ws.Range("E9").Value = '  -0.88%  '

$ws.Range("D10").Value = '19.48'
$ws.Range("E10").Value = '  -0.57%  '

$ws.Range("D11").Value = '0.0791'
$ws.Range("E11").Value = '  +1.04%  '

$ws.Range("D12").Value = '4.26'
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("D13").Value = '1.853.89'
$ws.Range("E13").Value = '  -0.20%  '

$ws.Range("D14").Value = '1.610.80'
$ws.Range("E14").Value = '  -1.70%  '

$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("E16").Value = '  -1.80%  '

$ws.Range("D17").Value = '62.79'
$ws.Range("E17").Value = '  -0.59%  '

$ws.Range("D18").Value = '25.723.29'
$ws.Range("E18").Value = '  -0.22%  '

$ws.Range("E19").Value = '  -0.33%  '

$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("D21").Value = '191.71'
$ws.Range("E21").Value = '  -1.06%  '

$ws.Range("D22").Value = '9.92'
$ws.Range("E22").Value = '  -0.21%  '

$ws.Range("E23").Value = '  +1.35%  '

$ws.Range("E24").Value = '  -0.50%  '

$ws.Range("E25").Value = '  +3.39%  '

$ws.Range("D26").Value = '141.69'
$ws.Range("E26").Value = '  +1.44%  '

$ws.Range("E27").Value = '  +1.96%  '

$ws.Range("E28").Value = '  -0.07%  '

$ws.Range("D29").Value = '15.46'
$ws.Range("E29").Value = '  -1.05%  '

$ws.Range("D30").Value = '1.24'
$ws.Range("E30").Value = '  -0.68%  '

$ws.Range("D31").Value = '0.0491'
$ws.Range("E31").Value = '  +0.57%  '

$ws.Range("E32").Value = '  -0.87%  '

$ws.Range("E33").Value = '  -1.30%  '

$ws.Range("E34").Value = '  -0.42%  '

$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").Value = '0.903'
$ws.Range("E36").Value = '  +0.76%  '

$ws.Range("D37").Value = '1.140.22'
$ws.Range("E37").Value = '  +3.22%  '

$ws.Range("E38").Value = '  -2.82%  '

$ws.Range("E39").Value = '  -1.95%  '

$ws.Range("E40").Value = '  -0.67%  '

$ws.Range("E41").Value = '  -0.16%  '

$ws.Range("E42").Value = '  -0.32%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '100.69'
$ws.Range("E43").Value = '  +1.39%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.54'
$ws.Range("E44").Value = '  -0.64%  '

$ws.Range("D45").Value = '0.805'
$ws.Range("E45").Value = '  -0.09%  '

$ws.Range("D46").Value = '1.763.15'
$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '55.16'
$ws.Range("E47").Value = '  +0.06%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.0509'
$ws.Range("E48").Value = '  +0.85%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '0.418'
$ws.Range("E49").Value = '  -0.30%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '1.44'
$ws.Range("E50").Value = '  +5.08%  '

$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = '2.35'
$ws.Range("E51").Value = '  -4.71%  '

$ws.Range("D2:D51").Style = "Normal"
